# TTrack_v1 sample_academic_transcript.xlsx
# v1.3.0 - Engine properly matches transcript x curriculum.
# Update subject codes/names to the real curriculum subjects, fix the
# Year value for the Elective row, and widen the Subject Name column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Software Engineering Principles
$ws.Range("C2").Value = "SEP401"
$ws.Range("D2").Value = "Software Engineering Principles"

# Row 3 - Software Development Management
$ws.Range("C3").Value = "SDM404"
$ws.Range("D3").Value = "Software Development Management"

# Row 4 - Research Methodologies
$ws.Range("C4").Value = "REM502"
$ws.Range("D4").Value = "Research Methodologies"

# Row 5 - Elective (also correct the Year back to 2024)
$ws.Range("A5").Value = 2024
$ws.Range("C5").Value = "Elective 1"
$ws.Range("D5").Value = "Elective Subject"

# Match the updated cells' font to Arial (no theme scheme) as applied in Excel
$ws.Range("C2:D4").Font.Name = "Arial"
$ws.Range("C5:D5").Font.Name = "Arial"

# Widen the Subject Name column to fit the longer subject names
$ws.Columns.Item(4).ColumnWidth = 26.8

"Edit complete"
